# Refresh the coin price / volume snapshot (GitHub Actions symbol-list update).
# Cells in columns D (Price) and E (Volume 1h) are numeric-looking text, so a
# leading apostrophe forces Excel to store them as text instead of
# auto-converting to Number/Percentage, matching the sheet's existing layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.43"
$ws.Range("E2").Value = "'6.27%"
$ws.Range("D3").Value = "'27.19"
$ws.Range("E3").Value = "'1.13%"
$ws.Range("E4").Value = "'1.47%"
$ws.Range("D5").Value = "'0.06275"
$ws.Range("E5").Value = "'0.95%"
$ws.Range("D6").Value = "'6.923"
$ws.Range("E6").Value = "'2.87%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8792"
$ws.Range("E7").Value = "'3.51%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.9436"
$ws.Range("E8").Value = "'3.33%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1459"
$ws.Range("E9").Value = "'3.91%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.05355"
$ws.Range("E10").Value = "'8.36%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07270"
$ws.Range("E11").Value = "'2.66%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03110"
$ws.Range("E12").Value = "'0.75%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09055"
$ws.Range("E13").Value = "'0.08%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001563"
$ws.Range("E14").Value = "'2.41%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006265"
$ws.Range("E15").Value = "'1.40%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005877"
$ws.Range("E16").Value = "'-1.46%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.446"
$ws.Range("E17").Value = "'0.00%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.264"
$ws.Range("E18").Value = "'2.92%"
$ws.Range("E19").Value = "'5.30%"
$ws.Range("D21").Value = "'0.1312"
$ws.Range("E21").Value = "'0.16%"
$ws.Range("D22").Value = "'3.858"
$ws.Range("E22").Value = "'-6.61%"
$ws.Range("D23").Value = "'0.04310"
$ws.Range("E23").Value = "'1.32%"
$ws.Range("D24").Value = "'0.001186"
$ws.Range("E24").Value = "'0.04%"
$ws.Range("D25").Value = "'0.004281"
$ws.Range("E25").Value = "'5.20%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("D27").Value = "'0.0001691"
$ws.Range("E27").Value = "'3.13%"
$ws.Range("D40").Value = "'0.04043"
$ws.Range("E40").Value = "'2.76%"
$ws.Range("D41").Value = "'0.006416"
$ws.Range("E41").Value = "'55.08%"
$ws.Range("D42").Value = "'0.1154"
$ws.Range("E42").Value = "'3.86%"
$ws.Range("D43").Value = "'0.002202"
$ws.Range("E43").Value = "'4.78%"
$ws.Range("D44").Value = "'0.01179"
$ws.Range("E44").Value = "'-11.46%"
$ws.Range("D45").Value = "'0.00005084"
$ws.Range("E45").Value = "'-1.55%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("D47").Value = "'2.378"
$ws.Range("E47").Value = "'843.75%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.02%"
